# stats_hrv_response_mist_phases.xlsx -- "added further analysis and result tables"
#
# The underlying analysis dropped the `test__effsize` parameter row (and its
# value `np2`) from the "parameter" sheet, and the repeated-measures ANOVA
# effect-size column was recomputed as generalized eta-squared ("ng2")
# instead of partial eta-squared ("np2"), with updated numeric results.

$wb = $excel.ActiveWorkbook

# --- 1. "parameter" sheet: remove the test__effsize / np2 row ----------
$wsParam = $wb.Worksheets.Item("parameter")
$wsParam.Rows.Item(7).Delete()

# --- 2. "rm_anova" sheet: rename effect-size column header np2 -> ng2 --
$wsAnova = $wb.Worksheets.Item("rm_anova")
$wsAnova.Cells.Item(2, 9).Value = "ng2"

# --- 3. "rm_anova" sheet: updated generalized eta-squared values -------
$wsAnova.Cells.Item(3, 9).Value = 0.0851
$wsAnova.Cells.Item(4, 9).Value = 0.0939
$wsAnova.Cells.Item(5, 9).Value = 0.2322
$wsAnova.Cells.Item(6, 9).Value = 0.0382
$wsAnova.Cells.Item(7, 9).Value = 0.0298
$wsAnova.Cells.Item(8, 9).Value = 0.1832
$wsAnova.Cells.Item(9, 9).Value = 0.0499
$wsAnova.Cells.Item(10, 9).Value = 0.0174
$wsAnova.Cells.Item(11, 9).Value = 0.1647
